# Update "想去人数" (want-to-go count) figures in column F on sheets
# "展览" and "全部类型" to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1083
$ws1.Range("F5").Value  = 409
$ws1.Range("F9").Value  = 6729
$ws1.Range("F10").Value = 150
$ws1.Range("F15").Value = 1080
$ws1.Range("F16").Value = 16081
$ws1.Range("F17").Value = 1577
$ws1.Range("F19").Value = 326
$ws1.Range("F20").Value = 174
$ws1.Range("F22").Value = 11294
$ws1.Range("F24").Value = 890
$ws1.Range("F25").Value = 4438
$ws1.Range("F26").Value = 296
$ws1.Range("F29").Value = 29
$ws1.Range("F30").Value = 316
$ws1.Range("F32").Value = 5216

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1083
$ws4.Range("F5").Value  = 409
$ws4.Range("F10").Value = 6729
$ws4.Range("F11").Value = 150
$ws4.Range("F17").Value = 1080
$ws4.Range("F18").Value = 16081
$ws4.Range("F19").Value = 1577
$ws4.Range("F21").Value = 326
$ws4.Range("F22").Value = 174
$ws4.Range("F26").Value = 11294
$ws4.Range("F28").Value = 890
$ws4.Range("F29").Value = 4438
$ws4.Range("F30").Value = 296
$ws4.Range("F33").Value = 29
$ws4.Range("F34").Value = 316
$ws4.Range("F36").Value = 5216
